$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: extend the data table down to row 409 by copying formatting
# (number formats / fonts) from the last existing row (382).
$ws.Range("A382:I382").Copy() | Out-Null
$ws.Range("A383:I409").PasteSpecial(-4122) | Out-Null

# Step 2: populate the new rows with the recorded training-log entries.
$newRows = @(
  @{ Row=383; A=45923; B="Amir Etien"; C=70; D=6; E=10; F=8; G="Pied et dos"; H=2 },
  @{ Row=384; A=45923; B="Yoann Martelat"; C=70; D=4; E=4; F=3; G="Genou"; H=7 },
  @{ Row=385; A=45923; B="Yoan Zouma"; C=70; D=4; E=9; F=6; G="Genou cheville pied"; H=2 },
  @{ Row=386; A=45923; B="Naim Ighbane"; C=70; D=6; E=6; F=4; G="Cheville droite"; H=7 },
  @{ Row=387; A=45923; B="Emmanuel Valey"; C=70; D=5; E=6; F=5; G="Ischio "; H=7 },
  @{ Row=388; A=45923; B="Hedi Nasri"; C=70; D=5; E=2; F=0; G=$null; H=8 },
  @{ Row=389; A=45923; B="Omar Benyounes"; C=70; D=5; E=3; F=0; G=$null; H=7 },
  @{ Row=390; A=45923; B="Wael Fareh"; C=70; D=7; E=1; F=2; G="Genou"; H=7 },
  @{ Row=391; A=45923; B="Kamal Bafounta"; C=70; D=6; E=5; F=4; G="Ischio"; H=5 },
  @{ Row=392; A=45923; B="Malik Boussaid"; C=70; D=2; E=1; F=0; G=$null; H=10 },
  @{ Row=393; A=45923; B="Jeremie Laurent"; C=70; D=7; E=6; F=0; G=$null; H=8 },
  @{ Row=394; A=45923; B="Karim Belmahi"; C=70; D=6; E=7; F=0; G=$null; H=10 },
  @{ Row=395; A=45923; B="Ilan Ihaddadene"; C=70; D=6; E=6; F=0; G=$null; H=8 },
  @{ Row=396; A=45923; B="Sofiane Belle"; C=70; D=4; E=4; F=2; G="Pied coup"; H=7 },
  @{ Row=397; A=45924; B="Kamal Bafounta"; C=70; D=6; E=5; F=5; G="Fesse cheville "; H=6 },
  @{ Row=398; A=45924; B="Naim Dhib"; C=70; D=5; E=5; F=4; G="Genou Tendon ischio"; H=7 },
  @{ Row=399; A=45924; B="Karim Belmahi"; C=70; D=8; E=7; F=1; G="Courbature"; H=10 },
  @{ Row=400; A=45924; B="Malik Boussaid"; C=70; D=2; E=0; F=0; G=$null; H=10 },
  @{ Row=401; A=45924; B="Omar Benyounes"; C=70; D=4; E=3; F=0; G=$null; H=7 },
  @{ Row=402; A=45924; B="Amir Etien"; C=70; D=7; E=9; F=5; G="Cheville coup"; H=7 },
  @{ Row=403; A=45924; B="Yoann Martelat"; C=70; D=6; E=5; F=3; G="Genou"; H=6 },
  @{ Row=404; A=45924; B="Yoan Zouma"; C=70; D=4; E=7; F=6; G="Pied"; H=2 },
  @{ Row=405; A=45924; B="Naim Ighbane"; C=70; D=5; E=6; F=4; G="Cheville droite"; H=7 },
  @{ Row=406; A=45924; B="Hedi Nasri"; C=70; D=5; E=6; F=2; G="Adducteur"; H=9 },
  @{ Row=407; A=45924; B="Ilan Ihaddadene"; C=70; D=7; E=7; F=0; G=$null; H=7 },
  @{ Row=408; A=45924; B="Emmanuel Valey"; C=70; D=7; E=7; F=5; G="Ischio"; H=7 },
  @{ Row=409; A=45924; B="Mattheo Haon"; C=70; D=6; E=8; F=0; G=$null; H=8 }
)

foreach ($row in $newRows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    if ($row.G -ne $null) {
        $ws.Cells.Item($r, 7).Value = $row.G
    }
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Formula = "=C" + $r + "*D" + $r
}

# Step 3: cells with a "Localisation douleur" entry must use the same
# text-font style (s=1) as other populated cells in column G, rather than
# the blank style (s=2) that was copied down in Step 1.
foreach ($row in $newRows) {
    if ($row.G -ne $null) {
        $r = $row.Row
        $ws.Range("G381").Copy() | Out-Null
        $ws.Cells.Item($r, 7).PasteSpecial(-4122) | Out-Null
    }
}

# Step 4: update the view state to match where the user ended up.
$ws.Range("K406").Select() | Out-Null
